$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells keep their original text representation
# (Price/Volume columns were stored as text strings, not numbers,
# so force Text number format before assigning the new values).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.221.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.969.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.019'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +1.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.54%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -6.90%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4180'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.58'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08835'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.095'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.328.80'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.92'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -7.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.863'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -8.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.385'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.020'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001098'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.50%  '

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -10.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06729'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.18'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -9.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.015'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.926'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.294.65'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.57%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.57'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.74'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.171'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -9.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.279'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -9.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.041'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -9.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09847'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.507'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -9.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.769'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.732'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02424'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -8.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.152'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -10.13%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06304'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6458'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.50'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -9.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2011'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -9.63%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6220'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -9.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.33'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -7.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.172'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.267'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -9.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.490'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000336'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06883'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.111'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -9.54%  '
